# ---------------------------------------------------------------------------
# Rebuild the "output" sheet (sheet2) into the new LaTeX-table layout, and
# update the shared strings / styles that go with it.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("output")

# Start from a clean sheet (clears values + formatting alike).
$ws.Cells.Clear()

# --- LaTeX table preamble -------------------------------------------------
$ws.Range("B2").Value = "\begin{table}[H]"
$ws.Range("B3").Value = "`t\begin{tabularx}{\columnwidth}"
$ws.Range("B4").Value = "{ >{\RaggedRight}p{9cm} | C }"
$ws.Range("B5").Value = "`t`t\hline"

# --- Header row (row 6) ----------------------------------------------------
$ws.Range("B6").Value = "Files migrated from Java to Kotlin (extension omitted)"
$ws.Range("C6").Value = "&"
$ws.Range("D6").Value = "Delta in lines of code"
$ws.Range("E6").Value = "\\"
$ws.Range("F6").Value = "\hline \hline"


# row 7: CameraRoll
$ws.Range("B7").Value = "CameraRoll"
$ws.Range("C7").Value = "&"
$ws.Range("D7").Formula = "=-data!O2+data!N3"
$ws.Range("E7").Value = "\\"
$ws.Range("F7").Value = "\hline"

# row 8: Classifier
$ws.Range("B8").Value = "Classifier"
$ws.Range("C8").Value = "&"
$ws.Range("D8").Formula = "=-data!O4+data!N5"
$ws.Range("E8").Value = "\\"
$ws.Range("F8").Value = "\hline"

# row 9: ListSingleton
$ws.Range("B9").Value = "ListSingleton"
$ws.Range("C9").Value = "&"
$ws.Range("D9").Formula = "=-data!O42+data!N43"
$ws.Range("E9").Value = "\\"
$ws.Range("F9").Value = "\hline"

# row 10: PermissionDenied
$ws.Range("B10").Value = "PermissionDenied"
$ws.Range("C10").Value = "&"
$ws.Range("D10").Formula = "=-data!O44+data!N45"
$ws.Range("E10").Value = "\\"
$ws.Range("F10").Value = "\hline"

# row 11: StartScreen
$ws.Range("B11").Value = "StartScreen"
$ws.Range("C11").Value = "&"
$ws.Range("D11").Formula = "=-data!O46+data!N47"
$ws.Range("E11").Value = "\\"
$ws.Range("F11").Value = "\hline"

# row 12: ViewFinder
$ws.Range("B12").Value = "ViewFinder"
$ws.Range("C12").Value = "&"
$ws.Range("D12").Formula = "=-data!O48+data!N49"
$ws.Range("E12").Value = "\\"
$ws.Range("F12").Value = "\hline"

# row 13: fragments/CameraRollPredictionsFragment
$ws.Range("B13").Value = "fragments/CameraRollPredictionsFragment"
$ws.Range("C13").Value = "&"
$ws.Range("D13").Formula = "=-data!O6+data!N7"
$ws.Range("E13").Value = "\\"
$ws.Range("F13").Value = "\hline"

# row 14: fragments/CameraSettingsFragment
$ws.Range("B14").Value = "fragments/CameraSettingsFragment"
$ws.Range("C14").Value = "&"
$ws.Range("D14").Formula = "=-data!O8+data!N9"
$ws.Range("E14").Value = "\\"
$ws.Range("F14").Value = "\hline"

# row 15: fragments/ModelSelectorFragment
$ws.Range("B15").Value = "fragments/ModelSelectorFragment"
$ws.Range("C15").Value = "&"
$ws.Range("D15").Formula = "=-data!O10+data!N11"
$ws.Range("E15").Value = "\\"
$ws.Range("F15").Value = "\hline"

# row 16: fragments/PredictionsFragment
$ws.Range("B16").Value = "fragments/PredictionsFragment"
$ws.Range("C16").Value = "&"
$ws.Range("D16").Formula = "=-data!O12+data!N13"
$ws.Range("E16").Value = "\\"
$ws.Range("F16").Value = "\hline"

# row 17: fragments/ProcessingUnitSelectorFragment
$ws.Range("B17").Value = "fragments/ProcessingUnitSelectorFragment"
$ws.Range("C17").Value = "&"
$ws.Range("D17").Formula = "=-data!O14+data!N15"
$ws.Range("E17").Value = "\\"
$ws.Range("F17").Value = "\hline"

# row 18: fragments/SmoothedPredictionsFragment
$ws.Range("B18").Value = "fragments/SmoothedPredictionsFragment"
$ws.Range("C18").Value = "&"
$ws.Range("D18").Formula = "=-data!O16+data!N17"
$ws.Range("E18").Value = "\\"
$ws.Range("F18").Value = "\hline"

# row 19: fragments/ThreadNumberFragment
$ws.Range("B19").Value = "fragments/ThreadNumberFragment"
$ws.Range("C19").Value = "&"
$ws.Range("D19").Formula = "=-data!O18+data!N19"
$ws.Range("E19").Value = "\\"
$ws.Range("F19").Value = "\hline"

# row 20: helpers/App
$ws.Range("B20").Value = "helpers/App"
$ws.Range("C20").Value = "&"
$ws.Range("D20").Formula = "=-data!O20+data!N21"
$ws.Range("E20").Value = "\\"
$ws.Range("F20").Value = "\hline"

# row 21: helpers/CameraEvents
$ws.Range("B21").Value = "helpers/CameraEvents"
$ws.Range("C21").Value = "&"
$ws.Range("D21").Formula = "=-data!O22+data!N23"
$ws.Range("E21").Value = "\\"
$ws.Range("F21").Value = "\hline"

# row 22: helpers/FreezeAnalyzer
$ws.Range("B22").Value = "helpers/FreezeAnalyzer"
$ws.Range("C22").Value = "&"
$ws.Range("D22").Formula = "=-data!O24+data!N25"
$ws.Range("E22").Value = "\\"
$ws.Range("F22").Value = "\hline"

# row 23: helpers/FreezeCallback
$ws.Range("B23").Value = "helpers/FreezeCallback"
$ws.Range("C23").Value = "&"
$ws.Range("D23").Formula = "=-data!O26+data!N27"
$ws.Range("E23").Value = "\\"
$ws.Range("F23").Value = "\hline"

# row 24: helpers/ImageUtils
$ws.Range("B24").Value = "helpers/ImageUtils"
$ws.Range("C24").Value = "&"
$ws.Range("D24").Formula = "=-data!O28+data!N29"
$ws.Range("E24").Value = "\\"
$ws.Range("F24").Value = "\hline"

# row 25: helpers/Logger
$ws.Range("B25").Value = "helpers/Logger"
$ws.Range("C25").Value = "&"
$ws.Range("D25").Formula = "=-data!O30+data!N31"
$ws.Range("E25").Value = "\\"
$ws.Range("F25").Value = "\hline"

# row 26: helpers/ModelConfig
$ws.Range("B26").Value = "helpers/ModelConfig"
$ws.Range("C26").Value = "&"
$ws.Range("D26").Formula = "=-data!O32+data!N33"
$ws.Range("E26").Value = "\\"
$ws.Range("F26").Value = "\hline"

# row 27: helpers/ProcessingUnit
$ws.Range("B27").Value = "helpers/ProcessingUnit"
$ws.Range("C27").Value = "&"
$ws.Range("D27").Formula = "=-data!O34+data!N35"
$ws.Range("E27").Value = "\\"
$ws.Range("F27").Value = "\hline"

# row 28: helpers/Recognition
$ws.Range("B28").Value = "helpers/Recognition"
$ws.Range("C28").Value = "&"
$ws.Range("D28").Formula = "=-data!O36+data!N37"
$ws.Range("E28").Value = "\\"
$ws.Range("F28").Value = "\hline"

# row 29: helpers/ResultItem
$ws.Range("B29").Value = "helpers/ResultItem"
$ws.Range("C29").Value = "&"
$ws.Range("D29").Formula = "=-data!O38+data!N39"
$ws.Range("E29").Value = "\\"
$ws.Range("F29").Value = "\hline"

# row 30: helpers/ResultItemComparator
$ws.Range("B30").Value = "helpers/ResultItemComparator"
$ws.Range("C30").Value = "&"
$ws.Range("D30").Formula = "=-data!O40+data!N41"
$ws.Range("E30").Value = "\\"
$ws.Range("F30").Value = "\hline"

# --- Totals row (row 31) ---------------------------------------------------
$ws.Range("B31").Value = "\bfseries{Cumulative Delta Over All Java Files}"
$ws.Range("B31").Font.Bold = $true
$ws.Range("C31").Value = "&"
$ws.Range("D31").Formula = '="\bfseries{"&SUM(D7:D30)&"}"'
$ws.Range("D31").Font.Bold = $true
$ws.Range("E31").Value = "\\"
$ws.Range("F31").Value = "\hline"

# --- LaTeX table closing ----------------------------------------------------
$ws.Range("B32").Value = "`t\end{tabularx}"
$ws.Range("B33").Value = "\caption[Java Files in TUM-Lens v1.0 And Change in Lines of Code After Their Conversion to Kotlin]{This table shows the results from the command line prompt in listing \ref{code:cloc}. Packages have been indicated as a prefix to the file name to resemble the original project structure of TUM-Lens v1.0. Overall, the total size of the codebase shrank because of the migration from Java to Kotlin. This is particularly impressive as understandably further logic needed to be added to existing classes in order to account for the new object detection functionality.}"
$ws.Range("B34").Value = "\label{tab:cloc}"
$ws.Range("B34").WrapText = $true
$ws.Rows.Item(34).RowHeight = 17
$ws.Range("B35").Value = "\end{table}"

# --- Column widths (character units; converted from the target pixel/char
#     widths the same way Excel stores ColumnWidth -> <col width=.../>) ----
$ws.Columns.Item(1).ColumnWidth = 1.6666666666666665   # -> 2.5
$ws.Columns.Item(2).ColumnWidth = 37.333333333333336    # -> ~38.1640625
$ws.Columns.Item(3).ColumnWidth = 1.6666666666666665   # -> 2.5
$ws.Columns.Item(4).ColumnWidth = 18.166666666666668    # -> 19

# --- Selection / view -------------------------------------------------------
$ws.Range("D6").Select()

Write-Output "done"
